$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

# ---------------------------------------------------------------------------
# 1) Column widths
#    - Column H gets wider (new, longer activity descriptions)
#    - Column G (TOTAL, previously unformatted) gets a plain "General" look
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 67.65          # -> stored width ~68.57 (H)
$ws.Columns.Item(7).ColumnWidth = 8.33           # -> stored width ~9.14 (G)

# ---------------------------------------------------------------------------
# 2) Build small "template" cells off to the side so we can stamp their
#    resulting styles onto many cells via Copy / PasteSpecial(xlPasteFormats).
#    This reproduces the formatting classes Excel created:
#      - plain cell, General number format (used for F/G data columns)
#      - centered text, no wrap (used for short ATIVIDADE entries)
#      - centered text, with wrap (used for the two long ATIVIDADE entries)
# ---------------------------------------------------------------------------
$tplNum = $ws.Range("J1")
$tplNum.NumberFormat = "General"

$tplCenter = $ws.Range("J2")
$tplCenter.HorizontalAlignment = $xlCenter

$tplWrap = $ws.Range("J3")
$tplWrap.HorizontalAlignment = $xlCenter
$tplWrap.WrapText = $true

# ---------------------------------------------------------------------------
# 3) Re-layout rows 4, 5 and 7: ENTRADA/SAÍDA values move from TURNO 01
#    (B/C) to TURNO 02 (D/E), leaving TURNO 01 blank but still time-formatted.
# ---------------------------------------------------------------------------
$ws.Range("B4:C4").Copy()
$ws.Range("D4:E4").PasteSpecial($xlPasteFormats)
$ws.Range("D4").Value = 0.58333333333333337
$ws.Range("E4").Value = 0.66666666666666663
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""

$ws.Range("B5:C5").Copy()
$ws.Range("D5:E5").PasteSpecial($xlPasteFormats)
$ws.Range("D5").Value = 0.625
$ws.Range("E5").Value = 0.70833333333333337
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = ""

$ws.Range("B7:C7").Copy()
$ws.Range("D7:E7").PasteSpecial($xlPasteFormats)
$ws.Range("D7").Value = 0.54166666666666663
$ws.Range("E7").Value = 0.70833333333333337
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""

# ---------------------------------------------------------------------------
# 4) Apply the plain General-format style to the existing F/G data cells
#    (rows 4-7) and give H4-H7 the new centered (no-wrap) text style.
# ---------------------------------------------------------------------------
$tplNum.Copy()
$ws.Range("F4:G7").PasteSpecial($xlPasteFormats)

$tplCenter.Copy()
$ws.Range("H4:H7").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 5) Append the new diary rows (8-17) plus the trailing totals row (18).
#    Each row only receives values in the columns the source row actually
#    used, so unused TURNO cells stay completely empty (no stray styling).
# ---------------------------------------------------------------------------

# Row 8: only TURNO 02 (D/E) has times; ht=45, wrapped activity text.
$ws.Range("A8").Value = 43810
$ws.Range("D8").Value = 0.54166666666666663
$ws.Range("E8").Value = 0.70833333333333337
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 15
$ws.Range("H8").Value = "Tentativa de instalação. Erro encontrado: tabela com tamanho grande. Solução: diminuir o tamanho do VARCHAR e alterar o tipo para TEXT quando for necessário."

# Row 9: both turnos used; ht=30, wrapped activity text.
$ws.Range("A9").Value = 43811
$ws.Range("B9").Value = 0.33333333333333331
$ws.Range("C9").Value = 0.41666666666666669
$ws.Range("D9").Value = 0.54166666666666663
$ws.Range("E9").Value = 0.70833333333333337
$ws.Range("F9").Value = 6
$ws.Range("G9").Value = 21
$ws.Range("H9").Value = "Instalação do java 8 (pc Renata) e instalação do sistema (foi colocado o xampp para iniciar automaticamente junto com o SO)"

# Row 10
$ws.Range("A10").Value = 43878
$ws.Range("B10").Value = 0.33333333333333331
$ws.Range("C10").Value = 0.41666666666666669
$ws.Range("D10").Value = 0.83333333333333337
$ws.Range("E10").Value = 0.91666666666666663
$ws.Range("F10").Value = 4
$ws.Range("G10").Value = 25
$ws.Range("H10").Value = "Consertando pdf (troca do iReport para o iText)"

# Row 11
$ws.Range("A11").Value = 43879
$ws.Range("B11").Value = 0.33333333333333331
$ws.Range("C11").Value = 0.45833333333333331
$ws.Range("D11").Value = 0.83333333333333337
$ws.Range("E11").Value = 0.95833333333333337
$ws.Range("F11").Value = 6
$ws.Range("G11").Value = 31
$ws.Range("H11").Value = "Consertando pdf (troca do iReport para o iText) e edição ícone (photoshop)"

# Row 12
$ws.Range("A12").Value = 43880
$ws.Range("B12").Value = 0.625
$ws.Range("C12").Value = 0.70833333333333337
$ws.Range("D12").Value = 0.83333333333333337
$ws.Range("E12").Value = 0.91666666666666663
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 35
$ws.Range("H12").Value = "Consertando pdf (troca do iReport para o iText)"

# Row 13: no TURNO 02 (D/E stay empty).
$ws.Range("A13").Value = 43881
$ws.Range("B13").Value = 0.54166666666666663
$ws.Range("C13").Value = 0.625
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 37
$ws.Range("H13").Value = "Reunião com o professor Leonardo e visita a enfermaria"

# Row 14: no TURNO 02 (D/E stay empty).
$ws.Range("A14").Value = 43883
$ws.Range("B14").Value = 0.625
$ws.Range("C14").Value = 0.79166666666666663
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 41
$ws.Range("H14").Value = "Edição do ícone (photoshop)"

# Row 15
$ws.Range("A15").Value = 43885
$ws.Range("B15").Value = 0.58333333333333337
$ws.Range("C15").Value = 0.75
$ws.Range("D15").Value = 0.83333333333333337
$ws.Range("E15").Value = 0.91666666666666663
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 47
$ws.Range("H15").Value = "PDF com o iText e finalização da edição do ícone (photoshop)"

# Row 16
$ws.Range("A16").Value = 43886
$ws.Range("B16").Value = 0.625
$ws.Range("C16").Value = 0.66666666666666663
$ws.Range("D16").Value = 0.83333333333333337
$ws.Range("E16").Value = 0.91666666666666663
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 50
$ws.Range("H16").Value = "PDF com o iText "

# Row 17
$ws.Range("A17").Value = 43887
$ws.Range("B17").Value = 0.4375
$ws.Range("C17").Value = 0.60416666666666663
$ws.Range("D17").Value = 0.83333333333333337
$ws.Range("E17").Value = 0.91666666666666663
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 56
$ws.Range("H17").Value = "PDF com o iText e mudança de ícone de cada frame"

# Row 18: running-total formula only, in column G.
$ws.Range("G18").Formula = "=F18+G17"

# ---------------------------------------------------------------------------
# 6) Formatting for the new rows, cell by cell so nothing beyond the real
#    used cells gets touched.
#    A8:A17 -> date format (copy from A4)
#    B/C/D/E used cells -> time format (copy from B4)
#    F/G used cells, plus G18 -> plain General format
#    H8,H9  -> centered + wrapped (long text)
#    H10:H17-> centered, no wrap (short text)
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A8:A17").PasteSpecial($xlPasteFormats)

$ws.Range("B4").Copy()
$ws.Range("D8:E9").PasteSpecial($xlPasteFormats)
$ws.Range("B9:E9").PasteSpecial($xlPasteFormats)
$ws.Range("B10:E12").PasteSpecial($xlPasteFormats)
$ws.Range("B13:C14").PasteSpecial($xlPasteFormats)
$ws.Range("B15:E17").PasteSpecial($xlPasteFormats)

$tplNum.Copy()
$ws.Range("F8:G17").PasteSpecial($xlPasteFormats)
$ws.Range("G18").PasteSpecial($xlPasteFormats)

$tplWrap.Copy()
$ws.Range("H8:H9").PasteSpecial($xlPasteFormats)

$tplCenter.Copy()
$ws.Range("H10:H17").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 7) Row heights for the two wrapped-text rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 30

# ---------------------------------------------------------------------------
# 8) Clean up the template scratch cells used for style stamping.
# ---------------------------------------------------------------------------
$ws.Range("J1:J3").Clear()

# ---------------------------------------------------------------------------
# 9) View state: selection moved to E21, mirroring the author's on-screen
#    state while editing (the sheet was scrolled so row 6 sits at the top).
# ---------------------------------------------------------------------------
$ws.Range("E21").Select()
